# TC05_INS_Filter_Doc-DCTD.xlsx -- "Fixed Bento 80 Test scripts"
#
# The Cypher query stored in startup!B2 gets an ORDER BY / LIMIT clause
# appended so the exported result set is capped and deterministically
# ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$cell = $ws.Range("B2")
$query = $cell.Value2
$cell.Value = $query + " ORDER BY p.project_id ASC LIMIT 100"

# The row holding the (now slightly different) wrapped query text is
# resized to match the re-flowed content.
$ws.Rows.Item(2).RowHeight = 244.8
